$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new mapping row (PROFIL / SO) was added above the existing "PROFIL / 0"
# row, pushing every row from the old row 11 down by one (old row 11 -> new
# row 12, ..., old row 18 -> new row 19).
$ws.Rows("11:11").Insert()

# Carry the formatting of the row immediately above (row 10) into the freshly
# inserted row 11, same as Excel does on a normal row insert.
$ws.Range("A10:P10").Copy()
$ws.Range("A11:P11").PasteSpecial(-4122)

# Fill in the new mapping entry: PROFIL / SO -> Spawalnia / Przygotowanie_do_spawania
$ws.Range("A11").Value = "PROFIL"
$ws.Range("B11").Value = "SO"
$ws.Range("C11").Value = "Spawalnia"
$ws.Range("D11").Value = "Przygotowanie_do_spawania"

# Leave the selection where the author left it after the edit
$ws.Range("B11").Select()
